$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (index 1) - counter (F column, "想去人数") bumps only
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 310
$ws1.Range("F5").Value = 8642
$ws1.Range("F7").Value = 10901
$ws1.Range("F20").Value = 416
$ws1.Range("F22").Value = 1845
$ws1.Range("F23").Value = 405
$ws1.Range("F24").Value = 594
$ws1.Range("F25").Value = 347
$ws1.Range("F27").Value = 71
$ws1.Range("F30").Value = 1231
$ws1.Range("F31").Value = 23
$ws1.Range("F32").Value = 7
$ws1.Range("F36").Value = 452
$ws1.Range("F42").Value = 360
$ws1.Range("F44").Value = 807
$ws1.Range("F45").Value = 649
$ws1.Range("F47").Value = 130
$ws1.Range("F48").Value = 119

# ---------------------------------------------------------------------
# Sheet "演出" (index 2) - a brand-new event was inserted as the new
# row 9, pushing the former rows 9-17 down to rows 10-18.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Push existing rows 9..17 down to 10..18, carrying values/formats along.
$ws2.Rows(9).Insert()

# The inserted row's A-cell lost its border/bold formatting; clone it
# from the (now shifted) row below so the whole A column stays uniform.
$ws2.Range("A10").Copy()
$ws2.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A is just the zero-based rank (row - 1); re-stamp the whole
# block so it reads 8,9,10,...,17 after the insert.
$ws2.Range("A9").Value = 8
$ws2.Range("A10").Value = 9
$ws2.Range("A11").Value = 10
$ws2.Range("A12").Value = 11
$ws2.Range("A13").Value = 12
$ws2.Range("A14").Value = 13
$ws2.Range("A15").Value = 14
$ws2.Range("A16").Value = 15
$ws2.Range("A17").Value = 16
$ws2.Range("A18").Value = 17

# New row 9 content - "北京·春日计划2024——特别二次元不插电音乐会"
# (B column holds plain-text dates, not real Excel dates, so force a
# text format before assigning or Excel will silently coerce the
# "2024-10-20" literal into a date serial.)
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "2024-10-20"
$ws2.Range("B9").Style = $ws2.Range("C9").Style
$ws2.Range("A9").Value = 8
$ws2.Range("C9").Value = "北京·春日计划2024——特别二次元不插电音乐会"
$ws2.Range("D9").Value = "复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)"
$ws2.Range("E9").Value = "2024.10.20 14:30-10.20 16:00"
$ws2.Range("F9").Value = 4
$ws2.Range("G9").Value = 119
$ws2.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=92853"
$ws2.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202409/scpLvBPg1727168336196.jpeg"

# The event that used to be last (old row 17, "花たん") picked up one
# more "want to go" vote now that it lives at row 18.
$ws2.Range("F18").Value = 391

# ---------------------------------------------------------------------
# Sheet "本地生活" (index 3) - counter bumps only
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 2826
$ws3.Range("F5").Value = 211

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4) - counter bumps, plus one title edit that
# dropped the "—【神秘嘉宾待官宣】" (mystery-guest-TBA) suffix.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value = 211
$ws4.Range("F9").Value = 8642
$ws4.Range("F11").Value = 10901
$ws4.Range("F19").Value = 416
$ws4.Range("F20").Value = 1845
$ws4.Range("F21").Value = 405
$ws4.Range("F22").Value = 594
$ws4.Range("F23").Value = 347
$ws4.Range("F25").Value = 71
$ws4.Range("F29").Value = 1231
$ws4.Range("F30").Value = 23
$ws4.Range("F36").Value = 452
$ws4.Range("C40").Value = "北京·法国姐姐”乔伊丝·乔纳森《小意思》巡回演唱会"
$ws4.Range("F41").Value = 360
$ws4.Range("F45").Value = 391
$ws4.Range("F46").Value = 649
$ws4.Range("F48").Value = 130
$ws4.Range("F49").Value = 119
